# Aggiunti segment per CIB - Corporate
# Adds two new "Business Line" rows (BL_CIB_000004 / BL_CIB_000005) to the
# "r Workspace_BusinessLine" sheet, and moves the active tab/selection back
# to the "Workspace" sheet (as it was when the file was last saved).

$wb = $excel.ActiveWorkbook

$wsBL = $wb.Worksheets.Item("r Workspace_BusinessLine")

# New row 6: WS_CIB_BL_CIB_000004 / BL_CIB_000004 (Business)
$wsBL.Range("A6").Value = "CREATE/MODIFY"
$wsBL.Range("B6").Value = "WS_CIB_BL_CIB_000004"
$wsBL.Range("C6").Value = "WS_CIB_BL_CIB_000004"
$wsBL.Range("E6").Value = "WS_CIB_BUSINESS"
$wsBL.Range("F6").Value = "BL_CIB_000004"

# New row 7: WS_CIB_BL_CIB_000005 / BL_CIB_000005 (Corporate)
$wsBL.Range("A7").Value = "CREATE/MODIFY"
$wsBL.Range("B7").Value = "WS_CIB_BL_CIB_000005"
$wsBL.Range("C7").Value = "WS_CIB_BL_CIB_000005"
$wsBL.Range("E7").Value = "WS_CIB_CORPORATE"
$wsBL.Range("F7").Value = "BL_CIB_000005"

# Leave the same cell selected on this sheet as in the saved file.
$wsBL.Range("E11").Select()

# Re-activate the "Workspace" sheet so it is the one marked as selected/active
# when the workbook is saved (matches the final state of the file).
$wsWorkspace = $wb.Worksheets.Item("Workspace")
$wsWorkspace.Activate()
